# Journal Lorenzo - add the "Kick-off du projet" sprint entries (rows 14-16)
# in the second weekly block of the journal, matching the author's
# "ajout code pour le robot" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Kick-off du projet - 1 day
$ws.Range("B14").Value = "Kick-off du projet"
$ws.Range("D14").Value = 1

# Row 15: Planification et revue des tâches - 0.5 day
$ws.Range("B15").Value = "Planification et revue des tâches"
$ws.Range("D15").Value = 0.5

# Row 16: Créaion de l'interface utilisateur - 1 day
$ws.Range("B16").Value = "Créaion de l'interface utilisateur"
$ws.Range("D16").Value = 1

# Reflect the cursor/selection ending on the last edited cell, like the author's session
$ws.Range("D16").Select()
